# MD_Double_Check.xlsx: add summary AVERAGE/ABS-deviation rows.
#
# "User 1" (sheet1) gets:
#   row 12  -> AVERAGE(col1:col10) per column
#   row 14  -> ABS(row1 - row12)/10 per column
#   row 15  -> ABS(row2 - row12)/10 per column
#   row 16  -> ABS(row3 - row12)/10 per column
# and becomes the active sheet/tab, with the final selection on A16:G16.
#
# "User 4" (sheet4) gets:
#   row 12  -> AVERAGE(col1:col10) per column
# and stops being the active sheet/tab, with the final selection on A12:G12.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("User 1")
$ws4 = $wb.Worksheets.Item("User 4")

# --- "User 1": averages row (12) ---
$ws1.Range("A12").Formula       = "=AVERAGE(A1:A10)"
$ws1.Range("B12:G12").Formula   = "=AVERAGE(B1:B10)"

# --- "User 1": absolute-deviation rows (14-16) ---
$ws1.Range("A14").Formula       = "=ABS(A1-A12)/10"
$ws1.Range("B14:G14").Formula   = "=ABS(B1-B12)/10"

$ws1.Range("A15").Formula       = "=ABS(A2-A12)/10"
$ws1.Range("B15:G15").Formula   = "=ABS(B2-B12)/10"

$ws1.Range("A16").Formula       = "=ABS(A3-A12)/10"
$ws1.Range("B16:G16").Formula   = "=ABS(B3-B12)/10"

# --- "User 4": averages row (12) ---
$ws4.Range("A12").Formula       = "=AVERAGE(A1:A10)"
$ws4.Range("B12:G12").Formula   = "=AVERAGE(B1:B10)"

# --- Final selections / active sheet: select "User 4" first, then
#     "User 1" last so "User 1" ends up the active tab (matches the
#     saved view state in the target workbook). ---
$null = $ws4.Range("A12:G12").Select()
$null = $ws1.Range("A16:G16").Select()
